$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-12: 45185 -> 45204
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value = 45204
    }
}
